# Add a new "time_taken" column (F) to the panel worksheet, matching the
# style used by the other header cells (B1:E1), and fill in a per-row
# timestamp for every data row (2-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell - same bold/border/centered style as the other headers
# (copy the formatting from E1 so the new column reuses the existing
# header cell style rather than creating a brand-new one).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$timestamps = @(
    "2021-10-05 10:52:16.436660",
    "2021-10-05 10:52:16.436670",
    "2021-10-05 10:52:16.436673",
    "2021-10-05 10:52:16.436676",
    "2021-10-05 10:52:16.436679",
    "2021-10-05 10:52:16.436681",
    "2021-10-05 10:52:16.436684",
    "2021-10-05 10:52:16.436686",
    "2021-10-05 10:52:16.436689",
    "2021-10-05 10:52:16.436692",
    "2021-10-05 10:52:16.436694",
    "2021-10-05 10:52:16.436697",
    "2021-10-05 10:52:16.436699",
    "2021-10-05 10:52:16.436702",
    "2021-10-05 10:52:16.436704",
    "2021-10-05 10:52:16.436707",
    "2021-10-05 10:52:16.436709",
    "2021-10-05 10:52:16.436712",
    "2021-10-05 10:52:16.436715",
    "2021-10-05 10:52:16.436717",
    "2021-10-05 10:52:16.436720",
    "2021-10-05 10:52:16.436722",
    "2021-10-05 10:52:16.436725",
    "2021-10-05 10:52:16.436728",
    "2021-10-05 10:52:16.436730",
    "2021-10-05 10:52:16.436733",
    "2021-10-05 10:52:16.436736",
    "2021-10-05 10:52:16.436738",
    "2021-10-05 10:52:16.436741",
    "2021-10-05 10:52:16.436743",
    "2021-10-05 10:52:16.436746",
    "2021-10-05 10:52:16.436749"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
